$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 340 (shifts rows 340:360 down to 343:363)
$ws.Rows("340:342").Insert()

# Row 340 (new) - Femacal de La Calera, Coquimbo, Ají, Americana (o), Primera
$ws.Cells.Item(340, 1).Value = 3
$ws.Cells.Item(340, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(340, 3).Value = "Coquimbo"
$ws.Cells.Item(340, 4).Value = 44585
$ws.Cells.Item(340, 5).Value = 5
$ws.Cells.Item(340, 6).Value = 100112021
$ws.Cells.Item(340, 7).Value = "Ají"
$ws.Cells.Item(340, 8).Value = "Americana (o)"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 73
$ws.Cells.Item(340, 11).Value = 14000
$ws.Cells.Item(340, 12).Value = 15000
$ws.Cells.Item(340, 13).Value = 14521
$ws.Cells.Item(340, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(340, 15).Value = "Limache"
$ws.Cells.Item(340, 16).Value = 968
$ws.Cells.Item(340, 17).Value = 15
$ws.Cells.Item(340, 18).Value = "Hortaliza"

# Row 341 (new) - Femacal de La Calera, Coquimbo, Ají, Americana (o), Primera
$ws.Cells.Item(341, 1).Value = 3
$ws.Cells.Item(341, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(341, 3).Value = "Coquimbo"
$ws.Cells.Item(341, 4).Value = 44585
$ws.Cells.Item(341, 5).Value = 5
$ws.Cells.Item(341, 6).Value = 100112021
$ws.Cells.Item(341, 7).Value = "Ají"
$ws.Cells.Item(341, 8).Value = "Americana (o)"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 73
$ws.Cells.Item(341, 11).Value = 23000
$ws.Cells.Item(341, 12).Value = 24000
$ws.Cells.Item(341, 13).Value = 23479
$ws.Cells.Item(341, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(341, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(341, 16).Value = 939
$ws.Cells.Item(341, 17).Value = 25
$ws.Cells.Item(341, 18).Value = "Hortaliza"

# Row 342 (new) - Femacal de La Calera, Coquimbo, Ají, Americana (o), Segunda
$ws.Cells.Item(342, 1).Value = 3
$ws.Cells.Item(342, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(342, 3).Value = "Coquimbo"
$ws.Cells.Item(342, 4).Value = 44585
$ws.Cells.Item(342, 5).Value = 5
$ws.Cells.Item(342, 6).Value = 100112021
$ws.Cells.Item(342, 7).Value = "Ají"
$ws.Cells.Item(342, 8).Value = "Americana (o)"
$ws.Cells.Item(342, 9).Value = "Segunda"
$ws.Cells.Item(342, 10).Value = 20
$ws.Cells.Item(342, 11).Value = 19000
$ws.Cells.Item(342, 12).Value = 19000
$ws.Cells.Item(342, 13).Value = 19000
$ws.Cells.Item(342, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(342, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(342, 16).Value = 760
$ws.Cells.Item(342, 17).Value = 25
$ws.Cells.Item(342, 18).Value = "Hortaliza"
